$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 'Critically low power Aircraft is landing .', 'Critically low power', '0-2', 'Event', 'Event'),
    @(0, 'Critically low power Aircraft is landing .', 'Aircraft is landing', '3-5', 'Event', 'Event'),
    @(1, 'Backward Obstacle Sensing is not functioning Ambient Light is too weak .', 'Ambient Light is too weak', '6-10', 'Event', 'Event'),
    @(1, 'Backward Obstacle Sensing is not functioning Ambient Light is too weak .', 'Backward Obstacle Sensing is not functioning', '0-5', 'Event', 'Event'),
    @(2, 'SD card speed low Change card .', 'Change card', '4-5', 'Event', 'Event'),
    @(2, 'SD card speed low Change card .', 'SD card speed low', '0-3', 'Event', 'Event'),
    @(3, 'Low Power Returning home .', 'Low Power', '0-1', 'Event', 'Event'),
    @(3, 'Low Power Returning home .', 'Returning home', '2-3', 'Event', 'Event'),
    @(4, 'SD card full Change card or delete files .', 'Change card or delete files', '3-7', 'Event', 'NonEvent'),
    @(4, 'SD card full Change card or delete files .', 'SD card full', '0-2', 'Event', 'Event'),
    @(5, 'Precision Landing Rectifying aircraft position .', 'Precision Landing', '0-1', 'Event', 'Event'),
    @(5, 'Precision Landing Rectifying aircraft position .', 'Rectifying aircraft position', '2-4', 'Event', 'Event'),
    @(6, 'Aircraft ActiveTrack available at max speed When exceeding nnn, Obstacle Avoidance is not available .', 'Aircraft ActiveTrack available at max speed', '0-5', 'NonEvent', 'NonEvent'),
    @(7, 'Mobile device CPU fully loaded Related performance will be affected .', 'Related performance will be affected', '5-9', 'Event', 'Event'),
    @(7, 'Mobile device CPU fully loaded Related performance will be affected .', 'Mobile device CPU fully loaded', '0-4', 'Event', 'Event'),
    @(8, 'Forward Obstacle Sensing is not functioning Ambient Light is too weak .', 'Ambient Light is too weak', '6-10', 'Event', 'Event'),
    @(8, 'Forward Obstacle Sensing is not functioning Ambient Light is too weak .', 'Forward Obstacle Sensing is not functioning', '0-5', 'Event', 'Event'),
    @(9, 'Extra payload detected Max altitude set to nnn and max distance set to nnn to ensure flight safety .', 'Extra payload detected', '0-2', 'Event', 'Event'),
    @(9, 'Extra payload detected Max altitude set to nnn and max distance set to nnn to ensure flight safety .', 'Max altitude set to nnn and max distance set to nnn to ensure flight safety', '3-17', 'Event', 'Event'),
    @(10, 'Battery error Aircraft returning to home automatically .', 'Battery error', '0-1', 'Event', 'Event'),
    @(10, 'Battery error Aircraft returning to home automatically .', 'Aircraft returning to home automatically', '2-6', 'Event', 'Event'),
    @(11, 'Precision Landing Locating Landing Point .', 'Precision Landing', '0-1', 'Event', 'Event'),
    @(11, 'Precision Landing Locating Landing Point .', 'Locating Landing Point', '2-4', 'Event', 'Event'),
    @(12, 'Current elevation above nnn Fly with caution .', 'Fly with caution', '4-6', 'NonEvent', 'NonEvent'),
    @(12, 'Current elevation above nnn Fly with caution .', 'Current elevation above nnn', '0-3', 'NonEvent', 'Event'),
    @(13, 'Aircraft is close to the Home Point Initiating Return to Home will now trigger Auto Landing .', 'Initiating Return to Home will now trigger Auto Landing', '7-15', 'Event', 'Event'),
    @(13, 'Aircraft is close to the Home Point Initiating Return to Home will now trigger Auto Landing .', 'Aircraft is close to the Home Point', '0-6', 'Event', 'NonEvent'),
    @(14, 'Forward vision sensor calibration error Auto calibration in progress .', 'Forward vision sensor calibration error', '0-4', 'Event', 'Event'),
    @(14, 'Forward vision sensor calibration error Auto calibration in progress .', 'Auto calibration in progress', '5-8', 'Event', 'Event'),
    @(15, 'Critically Low Voltage Aircraft will land .', 'Critically Low Voltage', '0-2', 'Event', 'Event'),
    @(15, 'Critically Low Voltage Aircraft will land .', 'Aircraft will land', '3-5', 'Event', 'Event'),
    @(16, 'Attitude is too large Backward Obstacle Sensing is not functioning .', 'Backward Obstacle Sensing is not functioning', '4-9', 'Event', 'Event'),
    @(16, 'Attitude is too large Backward Obstacle Sensing is not functioning .', 'Attitude is too large', '0-3', 'Event', 'Event'),
    @(17, 'Compass Error Compass Redundancy Switch .', 'Compass Error', '0-1', 'Event', 'Event'),
    @(17, 'Compass Error Compass Redundancy Switch .', 'Compass Redundancy Switch', '2-4', 'Event', 'Event'),
    @(18, 'Aircraft too Far Cannot enter Follow Me Mode .', 'Cannot enter Follow Me Mode', '3-7', 'Event', 'Event'),
    @(18, 'Aircraft too Far Cannot enter Follow Me Mode .', 'Aircraft too Far', '0-2', 'Event', 'Event'),
    @(19, 'Compass Error Compass direction is not the same with IMU .', 'Compass Error', '0-1', 'Event', 'Event'),
    @(19, 'Compass Error Compass direction is not the same with IMU .', 'Compass direction is not the same with IMU', '2-9', 'Event', 'Event'),
    @(20, 'Aircraft processor overheated Automatically returning to home .', 'Aircraft processor overheated', '0-2', 'Event', 'Event'),
    @(20, 'Aircraft processor overheated Automatically returning to home .', 'Automatically returning to home', '3-6', 'Event', 'Event'),
    @(21, 'Propeller Guard Mounted Propeller Guard mounted Forward Obstacle Sensing will be automatically switched off .', 'Forward Obstacle Sensing will be automatically switched off', '6-13', 'Event', 'Event'),
    @(22, 'Landing Protection Activated Aircraft will decelerate during landing .', 'Aircraft will decelerate during landing', '3-7', 'Event', 'Event'),
    @(22, 'Landing Protection Activated Aircraft will decelerate during landing .', 'Landing Protection Activated', '0-2', 'Event', 'Event'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
